$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "27.847.68"
$c.ClearFormats()
$ws.Range("E2").Value = "  -5.78%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.819.34"
$c.ClearFormats()
$ws.Range("E3").Value = "  -4.44%  "
$ws.Range("E4").Value = "  -0.43%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "328.14"
$c.ClearFormats()
$ws.Range("E5").Value = "  -2.91%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4625"
$c.ClearFormats()
$ws.Range("E7").Value = "  -2.96%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3839"
$c.ClearFormats()
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("E9").Value = "  -3.10%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.07837"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.65%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.9570"
$c.ClearFormats()
$ws.Range("E11").Value = "  -3.52%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "21.79"
$c.ClearFormats()
$ws.Range("E12").Value = "  -6.36%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "1.867.67"
$c.ClearFormats()
$ws.Range("E13").Value = "  -1.07%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.635"
$c.ClearFormats()
$ws.Range("E14").Value = "  -4.82%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.842"
$c.ClearFormats()
$ws.Range("E15").Value = "  -3.83%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.06851"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.46%  "
$ws.Range("E17").Value = "  -0.36%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "86.48"
$c.ClearFormats()
$ws.Range("E18").Value = "  -2.97%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.000009907"
$c.ClearFormats()
$ws.Range("E19").Value = "  -3.00%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "16.60"
$c.ClearFormats()
$ws.Range("E20").Value = "  -4.38%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E21").Value = "  -0.28%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "27.870.41"
$c.ClearFormats()
$ws.Range("E22").Value = "  -5.71%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.299"
$c.ClearFormats()
$ws.Range("E23").Value = "  -3.87%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "10.93"
$c.ClearFormats()
$ws.Range("E24").Value = "  -6.09%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.098"
$c.ClearFormats()
$ws.Range("E25").Value = "  -2.69%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.062.51"
$c.ClearFormats()
$ws.Range("E26").Value = "  -2.51%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "152.18"
$c.ClearFormats()
$ws.Range("E27").Value = "  -2.92%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.16"
$c.ClearFormats()
$ws.Range("E28").Value = "  -2.14%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.763"
$c.ClearFormats()
$ws.Range("E29").Value = "  -12.11%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "1.965"
$c.ClearFormats()
$ws.Range("E30").Value = "  -4.71%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "116.40"
$c.ClearFormats()
$ws.Range("E31").Value = "  -2.29%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.9316"
$c.ClearFormats()
$ws.Range("E32").Value = "  -6.83%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.09219"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.35%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.274"
$c.ClearFormats()
$ws.Range("E34").Value = "  -3.55%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.311"
$c.ClearFormats()
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.339"
$c.ClearFormats()
$ws.Range("E36").Value = "  -5.41%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.05922"
$c.ClearFormats()
$ws.Range("E37").Value = "  -8.56%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02141"
$c.ClearFormats()
$ws.Range("E38").Value = "  -4.48%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.140"
$c.ClearFormats()
$ws.Range("E39").Value = "  -4.44%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E40").Value = "  -0.25%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "7.537"
$c.ClearFormats()
$ws.Range("E41").Value = "  -2.52%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.5558"
$c.ClearFormats()
$ws.Range("E42").Value = "  -4.65%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "9.886"
$c.ClearFormats()
$ws.Range("E43").Value = "  -6.18%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.1762"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.23%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.215"
$c.ClearFormats()
$ws.Range("E45").Value = "  -3.96%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.212"
$c.ClearFormats()
$ws.Range("E46").Value = "  -10.02%  "
$ws.Range("E47").Value = "  -4.88%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.5229"
$c.ClearFormats()
$ws.Range("E48").Value = "  -4.58%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.06978"
$c.ClearFormats()
$ws.Range("E49").Value = "  -5.87%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.813"
$c.ClearFormats()
$ws.Range("E50").Value = "  -7.31%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "111.97"
$c.ClearFormats()
$ws.Range("E51").Value = "  -3.39%  "
